$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2:C10 ligand/receptor symbol columns (unaffected text, kept for completeness)
# and numeric columns E2:T10 with the recalculated TPM-derived values.

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.071327
$ws.Range("H2").Value = 0.213981
$ws.Range("I2").Value = 0.1064107741026886
$ws.Range("J2").Value = 0.1213845878939322
$ws.Range("M2").Value = 26.07194833333334
$ws.Range("N2").Value = 78.215845
$ws.Range("O2").Value = 0.7595928012803946
$ws.Range("P2").Value = 0.7818221335728009
$ws.Range("Q2").Value = 1.859633858771667
$ws.Range("R2").Value = 16.736704728945
$ws.Range("S2").Value = 0.08082885798707648
$ws.Range("T2").Value = 0.09490115749008927

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.071327
$ws.Range("H3").Value = 0.213981
$ws.Range("I3").Value = 0.1064107741026886
$ws.Range("J3").Value = 0.1213845878939322
$ws.Range("O3").Value = 0.155109029208254
$ws.Range("P3").Value = 0.1596482641062294
$ws.Range("Q3").Value = 0.379737672646
$ws.Range("R3").Value = 3.417639053814
$ws.Range("S3").Value = 0.01650527186836684
$ws.Range("T3").Value = 0.01937883874651631

# Row 4 (ECs -> MuSCs)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.071327
$ws.Range("H4").Value = 0.213981
$ws.Range("I4").Value = 0.1064107741026886
$ws.Range("J4").Value = 0.1213845878939322
$ws.Range("M4").Value = 2.927739
$ws.Range("N4").Value = 5.855478
$ws.Range("O4").Value = 0.08529816951135136
$ws.Range("P4").Value = 0.05852960232096958
$ws.Range("Q4").Value = 0.208826839653
$ws.Range("R4").Value = 1.252961037918
$ws.Range("S4").Value = 0.009076644247245247
$ws.Range("T4").Value = 0.007104591657326632

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.5235138962929958
$ws.Range("J5").Value = 0.5971812449832231
$ws.Range("M5").Value = 26.07194833333334
$ws.Range("N5").Value = 78.215845
$ws.Range("O5").Value = 0.7595928012803946
$ws.Range("P5").Value = 0.7818221335728009
$ws.Range("Q5").Value = 9.14892477094889
$ws.Range("R5").Value = 82.34032293854
$ws.Range("S5").Value = 0.3976573869944107
$ws.Range("T5").Value = 0.466889515082445

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.5235138962929958
$ws.Range("J6").Value = 0.5971812449832231
$ws.Range("O6").Value = 0.155109029208254
$ws.Range("P6").Value = 0.1596482641062294
$ws.Range("S6").Value = 0.08120173223103716
$ws.Range("T6").Value = 0.09533894911836849

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.5235138962929958
$ws.Range("J7").Value = 0.5971812449832231
$ws.Range("M7").Value = 2.927739
$ws.Range("N7").Value = 5.855478
$ws.Range("O7").Value = 0.08529816951135136
$ws.Range("P7").Value = 0.05852960232096958
$ws.Range("Q7").Value = 1.027374844316
$ws.Range("R7").Value = 6.164249065896
$ws.Range("S7").Value = 0.04465477706754797
$ws.Range("T7").Value = 0.03495278078240956

# Row 8 (MuSCs -> ECs)
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.248061
$ws.Range("H8").Value = 0.496122
$ws.Range("I8").Value = 0.3700753296043157
$ws.Range("J8").Value = 0.2814341671228447
$ws.Range("M8").Value = 26.07194833333334
$ws.Range("N8").Value = 78.215845
$ws.Range("O8").Value = 0.7595928012803946
$ws.Range("P8").Value = 0.7818221335728009
$ws.Range("Q8").Value = 6.467433575515001
$ws.Range("R8").Value = 38.80460145309
$ws.Range("S8").Value = 0.2811065562989075
$ws.Range("T8").Value = 0.2200314610002667

# Row 9 (MuSCs -> FAPs)
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.248061
$ws.Range("H9").Value = 0.496122
$ws.Range("I9").Value = 0.3700753296043157
$ws.Range("J9").Value = 0.2814341671228447
$ws.Range("O9").Value = 0.155109029208254
$ws.Range("P9").Value = 0.1596482641062294
$ws.Range("Q9").Value = 1.320651461778
$ws.Range("R9").Value = 7.923908770668
$ws.Range("S9").Value = 0.05740202510885004
$ws.Range("T9").Value = 0.04493047624134462

# Row 10 (MuSCs -> MuSCs)
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.248061
$ws.Range("H10").Value = 0.496122
$ws.Range("I10").Value = 0.3700753296043157
$ws.Range("J10").Value = 0.2814341671228447
$ws.Range("M10").Value = 2.927739
$ws.Range("N10").Value = 5.855478
$ws.Range("O10").Value = 0.08529816951135136
$ws.Range("P10").Value = 0.05852960232096958
$ws.Range("Q10").Value = 0.7262578640789999
$ws.Range("R10").Value = 2.905031456316
$ws.Range("S10").Value = 0.03156674819655815
$ws.Range("T10").Value = 0.0164722298812334

# Delete the rows for the removed "Resolving-Mac" sending cluster (former rows 11-13)
$ws.Range("A11:T13").EntireRow.Delete()
